# Applies the "OpponentWin" column addition to the Tabel1 table on Sheet1,
# updates the dependent SUMIF formulas in column J, and moves the active
# selection, matching the target commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# 1. Add a new table column ("OpponentWin") after the existing 5 columns.
$newCol = $tbl.ListColumns.Add()

# 2. Give the new header cell (F1) the same look as the other table header
#    cells (bold font + border, matching E1's style), then set its text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "OpponentWin"

# 3. Fill in the calculated column body (F2:F45) one cell at a time so each
#    cell carries its own (non-shared) formula, matching how Excel stores a
#    table calculated column.
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 6).Formula = "=IF(Tabel1[[#This Row],[Win rate]]<0.5,1,0)"
}

# 4. Update the SUMIF formulas in column J (rows 2-9) so the "Wins" total
#    counts wins both as Model and as Opponent.
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 10).Formula = "=SUMIF(Tabel1[Model],I$r,Tabel1[ModelWin])+SUMIF(Tabel1[Opponent],I$r,Tabel1[OpponentWin])"
}

# 5. Move the active selection to N13, matching the saved view state.
$ws.Range("N13").Select()
